# "added actors as descriptors"
# A new column "schlagworte_names" is inserted before the existing
# "schutzfrist" column (column U), shifting it and every later column
# one position to the right. The new column is then populated with the
# actor/descriptor names for the two publication rows that have them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at U; everything at/after old U shifts right.
$ws.Columns("U").Insert()

# New header for the inserted column.
$ws.Range("U1").Value = "schlagworte_names"

# New descriptor values for the two rows that have them.
$ws.Range("U5").Value = "Kränzle, Andreas; Meyerhans, Andreas"
$ws.Range("U6").Value = "Helg, Pater Lukas"

# Leave the selection where the edit finished.
$ws.Range("U7").Select()
